{"js": "// The document's body consists of a leading date paragraph followed by a\n// 20x5 table of arithmetic expressions (100 cells = 100 paragraphs), i.e.\n// 101 paragraphs total, each holding exactly one run/text. The commit\n// replaces the text of every one of those 101 paragraphs, in document\n// order, with a new value (new date, new expressions). Walking\n// `body.paragraphs` in order and rewriting each paragraph's text via\n// insertText(..., replace) preserves each run's existing formatting\n// (fonts, size, justification) while swapping the text content.\nconst newValues = [\"2023-05-29 Monday\", \"6-1=\", \"99-17=\", \"81-10=\", \"9+16=\", \"14+8=\", \"59-45=\", \"51+26=\", \"53-14=\", \"64+31=\", \"15+34=\", \"54-10=\", \"69-13=\", \"52-25=\", \"80+12=\", \"89-10=\", \"14+51=\", \"18-2=\", \"80+19=\", \"98-37=\", \"34-27=\", \"10+46=\", \"46+26=\", \"76+12=\", \"89-65=\", \"15+50=\", \"55-36=\", \"95-77=\", \"45-27=\", \"83-66=\", \"89-47=\", \"27-8=\", \"95-14=\", \"60-7=\", \"83-46=\", \"8+60=\", \"46+20=\", \"11+54=\", \"50+49=\", \"33-26=\", \"95-8=\", \"19+8=\", \"51-27=\", \"16+55=\", \"88+4=\", \"7+0=\", \"94-25=\", \"75-43=\", \"88+0=\", \"73+7=\", \"62-37=\", \"16+53=\", \"90-0=\", \"19+42=\", \"24+55=\", \"20+42=\", \"18+2=\", \"19+18=\", \"62-11=\", \"68-5=\", \"44+20=\", \"75+21=\", \"69-48=\", \"49+9=\", \"36+52=\", \"30-23=\", \"61-40=\", \"16+24=\", \"88-86=\", \"3+13=\", \"50-11=\", \"67-11=\", \"57+4=\", \"3+70=\", \"73-53=\", \"32+33=\", \"36+49=\", \"85-61=\", \"50+33=\", \"47+37=\", \"45-26=\", \"27-10=\", \"60+16=\", \"85-20=\", \"95-24=\", \"97-90=\", \"71-10=\", \"27+42=\", \"50+8=\", \"65+4=\", \"65-16=\", \"92-90=\", \"51+23=\", \"75-22=\", \"21+31=\", \"25+56=\", \"2+74=\", \"74+1=\", \"18-4=\", \"76+13=\", \"5+85=\"];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(\n    `Expected ${newValues.length} paragraphs but found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].insertText(newValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document's body consists of a leading centered date paragraph followed\n# by a single 20-row x 5-column table of arithmetic expressions (one\n# expression per cell). The commit rewrites the date and every cell's\n# expression. We address the date via the first paragraph and the cells via\n# the Tables collection (Cell(row, column)), which keeps each run's existing\n# character formatting (font/size) intact because only the Range.Text of an\n# already-formatted single run is being replaced.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday line at the top of the document.\n$d.Paragraphs.Item(1).Range.Text = \"2023-05-29 Monday\"\n\n# 2) Update every cell of the (only) table, row by row, left to right.\n$newValues = @(\n    @(\"6-1=\", \"99-17=\", \"81-10=\", \"9+16=\", \"14+8=\"),\n    @(\"59-45=\", \"51+26=\", \"53-14=\", \"64+31=\", \"15+34=\"),\n    @(\"54-10=\", \"69-13=\", \"52-25=\", \"80+12=\", \"89-10=\"),\n    @(\"14+51=\", \"18-2=\", \"80+19=\", \"98-37=\", \"34-27=\"),\n    @(\"10+46=\", \"46+26=\", \"76+12=\", \"89-65=\", \"15+50=\"),\n    @(\"55-36=\", \"95-77=\", \"45-27=\", \"83-66=\", \"89-47=\"),\n    @(\"27-8=\", \"95-14=\", \"60-7=\", \"83-46=\", \"8+60=\"),\n    @(\"46+20=\", \"11+54=\", \"50+49=\", \"33-26=\", \"95-8=\"),\n    @(\"19+8=\", \"51-27=\", \"16+55=\", \"88+4=\", \"7+0=\"),\n    @(\"94-25=\", \"75-43=\", \"88+0=\", \"73+7=\", \"62-37=\"),\n    @(\"16+53=\", \"90-0=\", \"19+42=\", \"24+55=\", \"20+42=\"),\n    @(\"18+2=\", \"19+18=\", \"62-11=\", \"68-5=\", \"44+20=\"),\n    @(\"75+21=\", \"69-48=\", \"49+9=\", \"36+52=\", \"30-23=\"),\n    @(\"61-40=\", \"16+24=\", \"88-86=\", \"3+13=\", \"50-11=\"),\n    @(\"67-11=\", \"57+4=\", \"3+70=\", \"73-53=\", \"32+33=\"),\n    @(\"36+49=\", \"85-61=\", \"50+33=\", \"47+37=\", \"45-26=\"),\n    @(\"27-10=\", \"60+16=\", \"85-20=\", \"95-24=\", \"97-90=\"),\n    @(\"71-10=\", \"27+42=\", \"50+8=\", \"65+4=\", \"65-16=\"),\n    @(\"92-90=\", \"51+23=\", \"75-22=\", \"21+31=\", \"25+56=\"),\n    @(\"2+74=\", \"74+1=\", \"18-4=\", \"76+13=\", \"5+85=\")\n)\n\n$table = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
